# bug fixes and code tuning
# - Flip a few "Yes"/"No" flag cells that were testing the wrong state.
# - Bump a couple of stale 2020 dates forward to 2021.
# - Nudge the remembered selection / window size left over from the last save.

$wb = $excel.ActiveWorkbook

# --- RoundTrip sheet -------------------------------------------------
$wsRoundTrip = $wb.Worksheets.Item("RoundTrip")

$wsRoundTrip.Range("I2").Value = "Yes"
$wsRoundTrip.Range("I3").Value = "No"
$wsRoundTrip.Range("E4").Value = "04/05/2021"
$wsRoundTrip.Range("I5").Value = "No"

# --- HotelSearch sheet -------------------------------------------------
$wsHotelSearch = $wb.Worksheets.Item("HotelSearch")

$wsHotelSearch.Range("F2").Value = "Yes"
$wsHotelSearch.Range("D3").Value = "05/04/2021"
$wsHotelSearch.Range("F3").Value = "Yes"

# --- restore the last-used selection on each sheet ----------------------
$wsRoundTrip.Activate()
$wsRoundTrip.Range("G13").Select()

$wsHotelSearch.Activate()
$wsHotelSearch.Range("D3").Select()

$wsRoundTrip.Activate()

Write-Host "edits applied"
